# Insert a new data row at row 224 (pushing the existing rows 224-242 down
# to 225-243), then populate the newly inserted row with the new record.
# This mirrors a weekly price-list update where a fresh entry is added at
# the top of this market's date-ordered block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 224.. down by one, inserting a blank row at 224.
$ws.Rows.Item(224).Insert()

# Populate the new row 224 with the new record's values.
$ws.Cells.Item(224, 1).Value  = 5
$ws.Cells.Item(224, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(224, 3).Value  = "Maule"
$ws.Cells.Item(224, 4).Value  = 45021
$ws.Cells.Item(224, 5).Value  = 7
$ws.Cells.Item(224, 6).Value  = 100112031
$ws.Cells.Item(224, 7).Value  = "Poroto verde"
$ws.Cells.Item(224, 8).Value  = "Sin especificar"
$ws.Cells.Item(224, 9).Value  = "Primera"
$ws.Cells.Item(224, 10).Value = 200
$ws.Cells.Item(224, 11).Value = 25000
$ws.Cells.Item(224, 12).Value = 25000
$ws.Cells.Item(224, 13).Value = 25000
$ws.Cells.Item(224, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(224, 15).Value = "Región del Maule"
$ws.Cells.Item(224, 16).Value = 1000
$ws.Cells.Item(224, 17).Value = 25
$ws.Cells.Item(224, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by
# the rest of column D.
$ws.Cells.Item(224, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
